# Automatische test-sync: 2025-08-03 18:36:50
#
# Adds the new "Testmail #11" row to the Logs sheet, adds the matching
# "Retour / Terugbetaling" tally row to the Dashboard sheet, extends the
# conditional formatting ranges and dimensions to cover the new row, and
# updates the Dashboard chart's category/value references accordingly.

$wb = $excel.ActiveWorkbook
$logsWs = $wb.Worksheets.Item("Logs")
$dashWs = $wb.Worksheets.Item("Dashboard")

# --- 1. Append new row 39 to the "Logs" sheet -----------------------------
$e39 = "Beste klant,`nDank u voor uw bericht. Kunt u ons meer informatie geven over uw retourzending, zoals uw ordernummer of retournummer, zodat we dit verder kunnen onderzoeken en u zo snel mogelijk kunnen helpen?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"

$logsWs.Range("A39").Value = "Mijn retour is nog steeds niet verwerkt."
$logsWs.Range("B39").Value = "mailmind.test@zohomail.eu"
$logsWs.Range("C39").Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$logsWs.Range("D39").Value = "Retour / Terugbetaling"
$logsWs.Range("E39").Value = $e39
$logsWs.Range("F39").Value = "2025-08-03 18:35:51"
$logsWs.Range("G39").Value = "Ja"
$logsWs.Range("H39").Value = "Nee"
$logsWs.Range("I39").Value = "Ja"
$logsWs.Range("J39").Value = "Nee"

# --- 2. Extend conditional formatting ranges on "Logs" from row 38 -> 39 --
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRng = $logsWs.Range($col + "2:" + $col + "38")
    $newRng = $logsWs.Range($col + "2:" + $col + "39")
    $fcs = $oldRng.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRng)
    }
}

# --- 3. Append new row 9 to the "Dashboard" sheet --------------------------
$dashWs.Range("A9").Value = "Retour / Terugbetaling"
$dashWs.Range("B9").Value = 1

# --- 4. Update chart series references from row 8 -> row 9 ----------------
$chartObj = $dashWs.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$9,Dashboard!`$B`$2:`$B`$9,1)"

Write-Host "Edit complete"
